{"js": "// Update the date line and every \"A\u00d7B=\" multiplication prompt in the table\n// to the new values, per the commit's regenerated problem set.\n\nconst body = context.document.body;\n\nconst replacements = [\n    [\"2024-07-22 Monday\", \"2024-07-23 Tuesday\"],\n    [\"94\u00d736=\", \"28\u00d746=\"],\n    [\"29\u00d717=\", \"65\u00d761=\"],\n    [\"85\u00d732=\", \"99\u00d798=\"],\n    [\"26\u00d793=\", \"82\u00d733=\"],\n    [\"28\u00d791=\", \"53\u00d727=\"],\n    [\"16\u00d798=\", \"29\u00d740=\"],\n    [\"62\u00d717=\", \"81\u00d761=\"],\n    [\"16\u00d735=\", \"44\u00d730=\"],\n    [\"91\u00d766=\", \"51\u00d719=\"],\n    [\"31\u00d792=\", \"44\u00d772=\"],\n    [\"84\u00d787=\", \"56\u00d762=\"],\n    [\"33\u00d725=\", \"37\u00d793=\"],\n    [\"80\u00d733=\", \"73\u00d765=\"],\n    [\"53\u00d762=\", \"65\u00d794=\"],\n    [\"25\u00d779=\", \"84\u00d723=\"],\n    [\"85\u00d759=\", \"65\u00d777=\"],\n    [\"34\u00d792=\", \"74\u00d721=\"],\n    [\"38\u00d757=\", \"30\u00d733=\"],\n    [\"31\u00d718=\", \"78\u00d756=\"],\n    [\"72\u00d773=\", \"21\u00d728=\"],\n    [\"74\u00d770=\", \"24\u00d742=\"],\n    [\"57\u00d718=\", \"15\u00d780=\"],\n    [\"86\u00d722=\", \"73\u00d711=\"],\n    [\"44\u00d740=\", \"73\u00d724=\"],\n    [\"17\u00d745=\", \"51\u00d749=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const r of results.items) {\n        r.insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the date line and every \"A\u00d7B=\" multiplication prompt in the table\n# to the new values, per the commit's regenerated problem set.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"2024-07-22 Monday\" = \"2024-07-23 Tuesday\"\n    \"94\u00d736=\"            = \"28\u00d746=\"\n    \"29\u00d717=\"            = \"65\u00d761=\"\n    \"85\u00d732=\"            = \"99\u00d798=\"\n    \"26\u00d793=\"            = \"82\u00d733=\"\n    \"28\u00d791=\"            = \"53\u00d727=\"\n    \"16\u00d798=\"            = \"29\u00d740=\"\n    \"62\u00d717=\"            = \"81\u00d761=\"\n    \"16\u00d735=\"            = \"44\u00d730=\"\n    \"91\u00d766=\"            = \"51\u00d719=\"\n    \"31\u00d792=\"            = \"44\u00d772=\"\n    \"84\u00d787=\"            = \"56\u00d762=\"\n    \"33\u00d725=\"            = \"37\u00d793=\"\n    \"80\u00d733=\"            = \"73\u00d765=\"\n    \"53\u00d762=\"            = \"65\u00d794=\"\n    \"25\u00d779=\"            = \"84\u00d723=\"\n    \"85\u00d759=\"            = \"65\u00d777=\"\n    \"34\u00d792=\"            = \"74\u00d721=\"\n    \"38\u00d757=\"            = \"30\u00d733=\"\n    \"31\u00d718=\"            = \"78\u00d756=\"\n    \"72\u00d773=\"            = \"21\u00d728=\"\n    \"74\u00d770=\"            = \"24\u00d742=\"\n    \"57\u00d718=\"            = \"15\u00d780=\"\n    \"86\u00d722=\"            = \"73\u00d711=\"\n    \"44\u00d740=\"            = \"73\u00d724=\"\n    \"17\u00d745=\"            = \"51\u00d749=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
